# exm: added pandoc lua-filter for dot in listing caption
#
# Adds a new paragraph style "ListingCaption" to the document's style
# sheet (styles.xml), based on Normal ("a"), quick-styled, no spacing
# after the paragraph, contextual spacing on, and bold run formatting.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$style = $d.Styles.Add("ListingCaption", 1)

# <w:basedOn w:val="a"/>  (styleId "a" == built-in "Normal")
$style.BaseStyle = "a"

# <w:qFormat/>
$style.QuickStyle = $true

# <w:pPr><w:spacing w:after="0"/> ... </w:pPr>
$style.ParagraphFormat.SpaceAfter = 0

# <w:pPr>...<w:contextualSpacing/></w:pPr>
$style.NoSpaceBetweenParagraphsOfSameStyle = $true

# <w:rPr><w:b/></w:rPr>
$style.Font.Bold = $true
